$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The Center of mass (CoG) is at the coordinate (2.410002574, " becomes
#    three runs: "The Center of mass (CoG) " / "of the whole object " /
#    "is at the coordinate (2.410002574, " (new middle clause inserted).
# ---------------------------------------------------------------------------
$piece1 = "The Center of mass (CoG) "
$piece2 = "of the whole object "
$piece3 = "is at the coordinate (2.410002574, "

$find = $d.Content
$find.Find.Execute("The Center of mass (CoG) is at the coordinate (2.410002574, ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ($piece1 + $piece2 + $piece3), 2) | Out-Null

$runStart = $find.Start
$boundary1 = $runStart + $piece1.Length
$boundary2 = $boundary1 + $piece2.Length

# Force the paragraph's run-coalescing pass to keep the two new boundaries by
# nudging (and then reverting) character formatting on the middle clause -
# this has to be the *last* edit touching the paragraph, otherwise later
# insertions re-merge every same-format run back together.
$mid = $d.Range($boundary1, $boundary2)
$mid.Bold = 1
$mid.Bold = 0

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the empty paragraph further down to sit
#    immediately after the ")" that now ends the CoG sentence.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$sentenceEnd = $d.Content
$sentenceEnd.Find.Execute("The Center of mass (CoG) of the whole object is at the coordinate (2.410002574, 5.327748134)") | Out-Null
$afterParen = $sentenceEnd.End

$target = $d.Range($afterParen, $afterParen)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
